$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3350

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 70
$ws.Range("I61").Value = 70
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 210
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -38
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 879.1548
$ws.Range("J129").Value = 967.5155999999999
$ws.Range("L129").Value = 2902.5468
$ws.Range("N129").Value = -12902.5468

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2489
$ws.Range("I141").Value = 2361.25
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 7083.75
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -1903.75
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5874.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6898.5713
$ws.Range("I88").Value = 2822.5
$ws.Range("J88").Value = 12333.333
$ws.Range("K88").Value = 2822.5
$ws.Range("L88").Value = 12333.333
$ws.Range("M88").Value = -2416.5
$ws.Range("N88").Value = -13145.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 6898.5713
$ws.Range("I91").Value = 2822.5
$ws.Range("J91").Value = 12333.333
$ws.Range("K91").Value = 2822.5
$ws.Range("L91").Value = 12333.333
$ws.Range("M91").Value = -1418.5
$ws.Range("N91").Value = -15141.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 37928.4
$ws.Range("J76").Value = 37928.4
$ws.Range("L76").Value = 37928.4
$ws.Range("N76").Value = -38558.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 37928.4
$ws.Range("J79").Value = 37928.4
$ws.Range("L79").Value = 37928.4
$ws.Range("N79").Value = -40112.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11881
$ws.Range("I86").Value = 13756.632
$ws.Range("J86").Value = 2971.75
$ws.Range("K86").Value = 13756.632
$ws.Range("L86").Value = 2971.75
$ws.Range("M86").Value = -12633.632
$ws.Range("N86").Value = -5217.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 11881
$ws.Range("I89").Value = 13756.632
$ws.Range("J89").Value = 2971.75
$ws.Range("K89").Value = 68783.16
$ws.Range("L89").Value = 14858.75
$ws.Range("M89").Value = -63167.16
$ws.Range("N89").Value = -26090.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1276.5385
$ws.Range("I94").Value = 600.25
$ws.Range("J94").Value = 1577.1111
$ws.Range("K94").Value = 600.25
$ws.Range("L94").Value = 1577.1111
$ws.Range("M94").Value = -149.25
$ws.Range("N94").Value = -2479.1111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1271.3334
$ws.Range("I134").Value = 1026.3636
$ws.Range("K134").Value = 3079.0908
$ws.Range("M134").Value = -544.0907999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2528.3809
$ws.Range("I31").Value = 1132.6538
$ws.Range("J31").Value = 4796.4375
$ws.Range("K31").Value = 1132.6538
$ws.Range("L31").Value = 4796.4375
$ws.Range("M31").Value = -837.6538
$ws.Range("N31").Value = -5386.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2528.3809
$ws.Range("I34").Value = 1132.6538
$ws.Range("J34").Value = 4796.4375
$ws.Range("K34").Value = 1132.6538
$ws.Range("L34").Value = 4796.4375
$ws.Range("M34").Value = -930.6538
$ws.Range("N34").Value = -5200.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3515
$ws.Range("I64").Value = 1140
$ws.Range("J64").Value = 3990
$ws.Range("K64").Value = 3420
$ws.Range("L64").Value = 11970
$ws.Range("M64").Value = -3150
$ws.Range("N64").Value = -12510

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3515
$ws.Range("I67").Value = 1140
$ws.Range("J67").Value = 3990
$ws.Range("K67").Value = 3420
$ws.Range("L67").Value = 11970
$ws.Range("M67").Value = -2484
$ws.Range("N67").Value = -13842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2881.2
$ws.Range("I70").Value = 1006
$ws.Range("J70").Value = 3350
$ws.Range("K70").Value = 3018
$ws.Range("L70").Value = 10050
$ws.Range("M70").Value = -2703
$ws.Range("N70").Value = -10680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2881.2
$ws.Range("I73").Value = 1006
$ws.Range("J73").Value = 3350
$ws.Range("K73").Value = 3018
$ws.Range("L73").Value = 10050
$ws.Range("M73").Value = -1926
$ws.Range("N73").Value = -12234

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1633.3334
$ws.Range("I75").Value = 450
$ws.Range("J75").Value = 4000
$ws.Range("K75").Value = 1350
$ws.Range("L75").Value = 12000
$ws.Range("M75").Value = -352
$ws.Range("N75").Value = -13996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1633.3334
$ws.Range("I78").Value = 450
$ws.Range("J78").Value = 4000
$ws.Range("K78").Value = 4050
$ws.Range("L78").Value = 36000
$ws.Range("M78").Value = 942
$ws.Range("N78").Value = -45984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 21061.875
$ws.Range("I87").Value = 9570.286
$ws.Range("J87").Value = 29999.777
$ws.Range("K87").Value = 28710.858
$ws.Range("L87").Value = 89999.33099999999
$ws.Range("M87").Value = -27462.858
$ws.Range("N87").Value = -92495.33099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 21061.875
$ws.Range("I90").Value = 9570.286
$ws.Range("J90").Value = 29999.777
$ws.Range("K90").Value = 86132.57399999999
$ws.Range("L90").Value = 269997.993
$ws.Range("M90").Value = -79892.57399999999
$ws.Range("N90").Value = -282477.993

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2679.9048
$ws.Range("I103").Value = 503.75
$ws.Range("J103").Value = 4019.077
$ws.Range("K103").Value = 1511.25
$ws.Range("L103").Value = 12057.231
$ws.Range("M103").Value = -632.25
$ws.Range("N103").Value = -13815.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1092.5714
$ws.Range("I122").Value = 627
$ws.Range("J122").Value = 1278.8
$ws.Range("K122").Value = 5643
$ws.Range("L122").Value = 11509.2
$ws.Range("M122").Value = -3193
$ws.Range("N122").Value = -16409.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2876264.2
$ws.Range("I129").Value = 1738.5714
$ws.Range("J129").Value = 5559155
$ws.Range("K129").Value = 5215.7142
$ws.Range("L129").Value = 16677465
$ws.Range("M129").Value = -215.7142000000003
$ws.Range("N129").Value = -16687465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1335.9333
$ws.Range("J131").Value = 1538.591
$ws.Range("L131").Value = 4615.772999999999
$ws.Range("N131").Value = -14695.773

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3717.44
$ws.Range("I134").Value = 1195.7333
$ws.Range("K134").Value = 3587.199900000001
$ws.Range("M134").Value = 1482.800099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 40140
$ws.Range("J100").Value = 40140
$ws.Range("L100").Value = 40140
$ws.Range("N100").Value = -42304

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 69580
$ws.Range("J101").Value = 69580
$ws.Range("L101").Value = 69580
$ws.Range("N101").Value = -76070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 27718.621
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 27718.621
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -37858.621

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 959.3
$ws.Range("I9").Value = 238.33333
$ws.Range("J9").Value = 2040.75
$ws.Range("K9").Value = 238.33333
$ws.Range("L9").Value = 2040.75
$ws.Range("M9").Value = -14.33332999999999
$ws.Range("N9").Value = -2488.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3474.75
$ws.Range("I16").Value = 1161.2307
$ws.Range("J16").Value = 13500
$ws.Range("K16").Value = 1161.2307
$ws.Range("L16").Value = 13500
$ws.Range("M16").Value = -991.2307000000001
$ws.Range("N16").Value = -13840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 37917.176
$ws.Range("J76").Value = 37917.176
$ws.Range("L76").Value = 37917.176
$ws.Range("N76").Value = -38593.176

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 37917.176
$ws.Range("J79").Value = 37917.176
$ws.Range("L79").Value = 37917.176
$ws.Range("N79").Value = -40257.176

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 45611.25
$ws.Range("J139").Value = 45611.25
$ws.Range("L139").Value = 45611.25
$ws.Range("N139").Value = -55891.25
